# cv124022a.xlsx — "correção nos dados e inicio da analise PNAD 2009"
#
# The "grupos de idade" sub-header row used a stray "unnamed: 1_level_1"
# label in B2 where it should read "total" (mirroring the B1 header).
# Also, two section-title rows ("situação do domicílio" and "grandes
# regiões e unidades da federação") were left as label-only rows with no
# data beneath them; they're removed so the data block is contiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled sub-header cell.
$ws.Range("B2").Value = "total"

# Remove the two empty section-header rows. Row 8 ("grandes regiões e
# unidades da federação") is deleted first so row 5's ("situação do
# domicílio") index doesn't shift before its own delete.
$ws.Rows("8").Delete()
$ws.Rows("5").Delete()
